# Leave Card update — 12/27/2023 4:01 PM upload
# Shifts the PERIOD end-dates in rows 94-107 of Sheet1 forward (extends the
# leave card by another month), fills in the newly-earned leave rows
# (103-105), and records a new "FL(1-0-0)" (Forced Leave) particular on
# row 105 together with its REMARKS date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Column A: PERIOD end-of-month dates shift forward one month ---
$ws.Range("A94").Value  = 44957
$ws.Range("A95").Value  = 44985
$ws.Range("A96").Value  = 45016
$ws.Range("A97").Value  = 45046
$ws.Range("A99").Value  = 45077
$ws.Range("A100").Value = 45107
$ws.Range("A101").Value = 45138
$ws.Range("A102").Value = 45169
$ws.Range("A103").Value = 45199
$ws.Range("A104").Value = 45230
$ws.Range("A105").Value = 45260
$ws.Range("A106").Value = 45291
$ws.Range("A107").Value = 45322

# --- Newly earned leave credits for the three added periods ---
$ws.Range("C103").Value = 1.25
$ws.Range("C104").Value = 1.25

# Row 105 — new "Forced Leave" entry: particulars, earned credit, the day
# used, and the remarks date. Copy K103's number format (date) onto K105
# first so the new date value picks up the same style Excel would apply.
$ws.Range("K103").Copy()
$ws.Range("K105").PasteSpecial(-4122)

$ws.Range("B105").Value = "FL(1-0-0)"
$ws.Range("C105").Value = 1.25
$ws.Range("D105").Value = 1
$ws.Range("K105").Value = 45258

# --- Restore the on-screen selection to where the editor left off ---
$ws.Range("F113").Select()
